$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Conrol" typo -> "Control"
$ws.Range("B1").Value = "Control"

# Copy the header style (bold white-on-blue, style index used by A1/B1) onto
# the new header cells C1:G1 before setting their text, so they pick up the
# same cellXf instead of Excel minting a near-duplicate style.
$ws.Range("A1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C1").Value = "Descargar Acuse"
$ws.Range("C2").Value = "SI"
$ws.Range("G1").Value = "Fila"
$ws.Range("E1").Value = "Directorio de Guardado"
$ws.Range("F1").Value = "Nombre de Guardado"
$ws.Range("D1").Value = "ST"

# Fill in the SI defaults for the new C/D columns, rows 2-5
$ws.Range("C2:C5").Value = "SI"
$ws.Range("D2:D5").Value = "SI"

# Row-number helper column
$ws.Range("G2").Formula = "=ROW()"
$ws.Range("G3").Formula = "=ROW()"
$ws.Range("G4").Formula = "=ROW()"
$ws.Range("G5").Formula = "=ROW()"

# Data validation dropdown (SI/NO) on C2:D5
$validationRange = $ws.Range("C2:D5")
$validationRange.Validation.Add(3, 1, 1, '"SI,NO"')
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

# Widen column A (engine quantizes ColumnWidth to 1/6-character steps, so
# 11.8 is the closest input that lands on the saved file's stored width)
$ws.Columns.Item(1).ColumnWidth = 11.8

# Move the active selection to D5, matching the saved cursor position
$ws.Range("D5").Select() | Out-Null
